# Update "想去人数" (number of people interested) figures for several events
# on the "展览" and "全部类型" sheets, as reflected in the regenerated data
# snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1): rows 2-6, column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1044
$wsExhibit.Range("F3").Value = 247
$wsExhibit.Range("F4").Value = 2622
$wsExhibit.Range("F5").Value = 45
$wsExhibit.Range("F6").Value = 561

# Sheet "全部类型" (sheetId 4): rows 4-8, column F (same events, offset by 2 rows)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1044
$wsAll.Range("F5").Value = 247
$wsAll.Range("F6").Value = 2622
$wsAll.Range("F7").Value = 45
$wsAll.Range("F8").Value = 561
